# Apply the updates described by the commit "actualizacion graficos y tablas"
# 1. Add 2019 and 2020 rows of data to the "Datos" sheet
# 2. Update the definition text (row with DEFINICIÓN) on "Ficha técnica" sheet
# 3. Update the CITA text on "Ficha técnica" sheet

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Datos" ---
$wsDatos = $wb.Worksheets.Item("Datos")

# The "Fecha" (year) column stores its values as text (matching the existing
# 2012/2014/.../2018 cells), so force text formatting before writing the
# numeric-looking year labels to keep them from being auto-converted to numbers.
$wsDatos.Range("A8:A9").NumberFormat = "@"

# New row for 2019
$wsDatos.Range("A8").Value = "2019"
$wsDatos.Range("B8").Value = 5.2
$wsDatos.Range("C8").Value = 2.3
$wsDatos.Range("D8").Value = 5.9
$wsDatos.Range("E8").Value = 4.1

# New row for 2020
$wsDatos.Range("A9").Value = "2020"
$wsDatos.Range("B9").Value = 4.7
$wsDatos.Range("C9").Value = 2
$wsDatos.Range("D9").Value = 6.4
$wsDatos.Range("E9").Value = 3.7

# --- Sheet 2: "Ficha técnica" ---
$wsFicha = $wb.Worksheets.Item("Ficha técnica")

# Update DEFINICIÓN value (B5, next to A5 = "DEFINICIÓN")
$wsFicha.Range("B5").Value = "El indicador mide la razón de las consultas médicas no urgentes sobre las consultas médicas urgentes y busca reflejar el control periódico, la prevención y detección temprana de las afecciones de salud. El Total SNIS es el promedio ponderado por cantidad de afiliados de cada tipo de prestador."

# Update CITA value (B8, next to A8 = "CITA")
$wsFicha.Range("B8").Value = "UMAD con base en SINADI - MSP (*2020 datos preliminares)"
